$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '37.791.55'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  +1.12%  '
$c = $ws.Range('D3')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '2.085.46'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('E4').Value = '  +0.05%  '
$c = $ws.Range('D5')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '234.45'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('E6').Value = '  +0.37%  '
$c = $ws.Range('D7')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '58.83'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +3.06%  '
$ws.Range('E8').Value = '  -0.02%  '
$c = $ws.Range('D9')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '0.392'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +0.12%  '
$c = $ws.Range('D10')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '0.0791'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +1.90%  '
$ws.Range('E11').Value = '  +2.72%  '
$c = $ws.Range('D12')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '2.395.15'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +0.95%  '
$c = $ws.Range('D13')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '14.74'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +2.84%  '
$c = $ws.Range('D14')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '21.19'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +2.72%  '
$c = $ws.Range('D15')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '0.768'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -1.15%  '
$ws.Range('E16').Value = '  +2.19%  '
$c = $ws.Range('D17')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '2.090.01'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +1.07%  '
$c = $ws.Range('D18')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '37.747.52'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +1.11%  '
$c = $ws.Range('D19')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '6.19'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -0.64%  '
$c = $ws.Range('D20')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '71.26'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +2.50%  '
$c = $ws.Range('D21')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '0.0₃0839'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +2.48%  '
$c = $ws.Range('D22')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '228.81'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +1.09%  '
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('E24').Value = '  -1.02%  '
$c = $ws.Range('D25')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '2.38'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -0.27%  '
$c = $ws.Range('D26')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '170.21'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +1.41%  '
$ws.Range('E27').Value = '  +5.85%  '
$c = $ws.Range('D28')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '9.00'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +1.71%  '
$ws.Range('E29').Value = '  +0.01%  '
$c = $ws.Range('D30')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '19.52'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +2.30%  '
$ws.Range('E31').Value = '  +2.56%  '
$c = $ws.Range('D32')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '4.69'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +3.06%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D33')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '0.0632'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +2.45%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D34')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '4.71'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +3.86%  '
$c = $ws.Range('D35')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '2.51'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +1.77%  '
$c = $ws.Range('D36')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '3.46'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +2.75%  '
$c = $ws.Range('D37')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '1.83'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +3.14%  '
$ws.Range('E38').Value = '  +0.04%  '
$c = $ws.Range('D39')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '5.39'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -3.96%  '
$c = $ws.Range('D40')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '0.0990'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +4.06%  '
$ws.Range('E41').Value = '  +0.43%  '
$c = $ws.Range('D42')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '98.71'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +1.89%  '
$ws.Range('B43').Value = 'FTXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$c = $ws.Range('D43')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '4.39'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +4.10%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D44')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '0.0215'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +1.21%  '
$c = $ws.Range('D45')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '1.465.29'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -1.66%  '
$ws.Range('E46').Value = '  +0.62%  '
$ws.Range('E47').Value = '  +4.01%  '
$c = $ws.Range('D48')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '16.01'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +5.23%  '
$ws.Range('E49').Value = '  +2.63%  '
$ws.Range('E50').Value = '  +2.37%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$c = $ws.Range('D51')
$c.Style = 'Normal'
$c.NumberFormat = '@'
$c.Value = '2.279.75'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +0.87%  '
